$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Final (post-edit) resource table, rows 2-16 (columns A-D)
$data = @(
    @("vgr_2", 885,  1, "/vgr/pick_up_and_transport"),
    @("hbw_2", 1581, 2, "/hbw/store_empty_bucket, /hbw/unload"),
    @("mm_2",  381,  4, "/mm/deburr, /mm/drill, /mm/mill, /mm/transport_from_to"),
    @("pm_1",  204,  3, "/pm/punch_gill, /pm/punch_recesses, /pm/punch_ribbing"),
    @("ov_2",  330,  1, "/ov/burn"),
    @("hbw_1", 873,  4, "/hbw/get_empty_bucket, /hbw/store, /hbw/store_empty_bucket, /hbw/unload"),
    @("ov_1",  612,  2, "/ov/burn, /ov/temper"),
    @("sm_2",  309,  2, "/sm/sort, /sm/transport"),
    @("vgr_1", 1866, 1, "/vgr/pick_up_and_transport"),
    @("wt_1",  447,  1, "/wt/pick_up_and_transport"),
    @("sm_1",  378,  2, "/sm/sort, /sm/transport"),
    @("mm_1",  576,  4, "/mm/deburr, /mm/drill, /mm/mill, /mm/transport_from_to"),
    @("dm_2",  177,  3, "/dm/cylindrical_drill, /dm/drill, /dm/lower"),
    @("hw_1",  522,  1, "/hw/human_review"),
    @("wt_2",  330,  1, "/wt/pick_up_and_transport")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
}
